$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "sheet 01"

# New column widths for O (15) and P (16) -- closest achievable quantized values
$ws.Columns.Item(15).ColumnWidth = 12.6
$ws.Columns.Item(16).ColumnWidth = 14.45

# Header row additions
$ws.Range("O1").Value = "degree"
$ws.Range("P1").Value = "isOpentoWork"

# Per-row data updates
# Row 2
$ws.Range("A2").Value = 17001
$ws.Range("F2").Value = 2017
$ws.Range("H2").Value = $false
$ws.Range("O2").Value = "CS"
$ws.Range("P2").Value = $false

# Row 3
$ws.Range("A3").Value = 18002
$ws.Range("F3").Value = 2018
$ws.Range("H3").Value = $false
$ws.Range("O3").Value = "CS"
$ws.Range("P3").Value = $false

# Row 4
$ws.Range("A4").Value = 18003
$ws.Range("F4").Value = 2018
$ws.Range("H4").Value = $false
$ws.Range("O4").Value = "CS"
$ws.Range("P4").Value = $false

# Row 5
$ws.Range("A5").Value = 17004
$ws.Range("F5").Value = 2017
$ws.Range("H5").Value = $false
$ws.Range("O5").Value = "CS"
$ws.Range("P5").Value = $false

# Row 6
$ws.Range("A6").Value = 17005
$ws.Range("F6").Value = 2017
$ws.Range("H6").Value = $false
$ws.Range("O6").Value = "CS"
$ws.Range("P6").Value = $false

# Row 7
$ws.Range("A7").Value = 17006
$ws.Range("F7").Value = 2017
$ws.Range("H7").Value = $false
$ws.Range("O7").Value = "CS"
$ws.Range("P7").Value = $false

# Row 8
$ws.Range("A8").Value = 17317
$ws.Range("F8").Value = 2017
$ws.Range("H8").Value = $false
$ws.Range("O8").Value = "DS"
$ws.Range("P8").Value = $false

# Row 9
$ws.Range("A9").Value = 17318
$ws.Range("F9").Value = 2017
$ws.Range("H9").Value = $false
$ws.Range("O9").Value = "DS"
$ws.Range("P9").Value = $false

# Row 10
$ws.Range("A10").Value = 17319
$ws.Range("F10").Value = 2017
$ws.Range("H10").Value = $false
$ws.Range("O10").Value = "DS"
$ws.Range("P10").Value = $false

# Row 11
$ws.Range("A11").Value = 18320
$ws.Range("F11").Value = 2018
$ws.Range("H11").Value = $false
$ws.Range("O11").Value = "DS"
$ws.Range("P11").Value = $false

# Row 12
$ws.Range("A12").Value = 18411
$ws.Range("F12").Value = 2018
$ws.Range("H12").Value = $false
$ws.Range("O12").Value = "DS"
$ws.Range("P12").Value = $false

# Row 13
$ws.Range("A13").Value = 18412
$ws.Range("F13").Value = 2018
$ws.Range("H13").Value = $false
$ws.Range("O13").Value = "DS"
$ws.Range("P13").Value = $false

# Row 14
$ws.Range("A14").Value = 17413
$ws.Range("F14").Value = 2017
$ws.Range("H14").Value = $false
$ws.Range("O14").Value = "DS"
$ws.Range("P14").Value = $false

# Row 15
$ws.Range("A15").Value = 17414
$ws.Range("F15").Value = 2017
$ws.Range("H15").Value = $false
$ws.Range("O15").Value = "DS"
$ws.Range("P15").Value = $false

# Row 16
$ws.Range("A16").Value = 17415
$ws.Range("F16").Value = 2017
$ws.Range("H16").Value = $false
$ws.Range("O16").Value = "CS"
$ws.Range("P16").Value = $false

# Row 17
$ws.Range("A17").Value = 17516
$ws.Range("F17").Value = 2017
$ws.Range("H17").Value = $false
$ws.Range("O17").Value = "DS"
$ws.Range("P17").Value = $false

# Row 18
$ws.Range("A18").Value = 17517
$ws.Range("F18").Value = 2017
$ws.Range("H18").Value = $false
$ws.Range("O18").Value = "STAT"
$ws.Range("P18").Value = $false

# Row 19
$ws.Range("A19").Value = 18518
$ws.Range("F19").Value = 2018
$ws.Range("H19").Value = $false
$ws.Range("O19").Value = "STAT"
$ws.Range("P19").Value = $false

# Row 20
$ws.Range("A20").Value = 18519
$ws.Range("F20").Value = 2018
$ws.Range("H20").Value = $false
$ws.Range("O20").Value = "STAT"
$ws.Range("P20").Value = $false

# Row 21
$ws.Range("A21").Value = 18520
$ws.Range("F21").Value = 2018
$ws.Range("H21").Value = $false
$ws.Range("O21").Value = "STAT"
$ws.Range("P21").Value = $false

# Row 22
$ws.Range("A22").Value = 17821
$ws.Range("F22").Value = 2017
$ws.Range("H22").Value = $false
$ws.Range("O22").Value = "STAT"
$ws.Range("P22").Value = $false

# Row 23
$ws.Range("A23").Value = 17822
$ws.Range("F23").Value = 2017
$ws.Range("H23").Value = $false
$ws.Range("O23").Value = "STAT"
$ws.Range("P23").Value = $false

# Row 24
$ws.Range("A24").Value = 17823
$ws.Range("F24").Value = 2017
$ws.Range("H24").Value = $false
$ws.Range("O24").Value = "STAT"
$ws.Range("P24").Value = $false

# Row 25
$ws.Range("A25").Value = 17824
$ws.Range("F25").Value = 2017
$ws.Range("H25").Value = $false
$ws.Range("O25").Value = "STAT"
$ws.Range("P25").Value = $false

# Row 26
$ws.Range("A26").Value = 17825
$ws.Range("F26").Value = 2017
$ws.Range("H26").Value = $false
$ws.Range("O26").Value = "STAT"
$ws.Range("P26").Value = $false

# Selection
$ws.Range("N28").Select()
